# Generate Report for Archive
# - Flip the "Ready for handoff" status to "In Translation" everywhere it
#   appears (Overview summary columns for zh-cn/de-de, plus each language
#   sheet's own Status column) so every reference collapses back onto a
#   single shared string.
# - Re-size the now-narrower Status columns to match the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column C, rows 2-3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column C, rows 2-3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Shrink the Status columns now that the text is shorter ---
# (closest width the host's ColumnWidth setter can represent)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
